$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Book_title"
$ws.Range("B1").Value = "Author"
$ws.Range("C1").Value = "ISBN"
$ws.Range("D1").Value = "Genre"
$ws.Range("E1").Value = "Availability_status"

$ws.Range("A2").Value = "Harry_potter"
$ws.Range("B2").Value = "maaz"
$ws.Range("C2").Value = "123-sdfg-456"
$ws.Range("D2").Value = "Horror"
$ws.Range("E2").Value = 0

$ws.Range("A3").Value = "Ironman"
$ws.Range("B3").Value = "tony"
$ws.Range("C3").Value = "3456fg23-32ss-2"
$ws.Range("D3").Value = "Action"
$ws.Range("E3").Value = 29

$ws.Range("A4").Value = "Gray_Hat_hakking"
$ws.Range("B4").Value = "Allen_Harper"
$ws.Range("C4").Value = "0-07-107731-6"
$ws.Range("D4").Value = "Education"
$ws.Range("E4").Value = 9

$ws.Range("A5").Value = "Hacking"
$ws.Range("B5").Value = "umar"
$ws.Range("C5").Value = "1234567890-z"
$ws.Range("D5").Value = "Education"
$ws.Range("E5").Value = 20

$ws.Range("A6").Value = "Java_book"
$ws.Range("B6").Value = "umar_khan"
$ws.Range("C6").Value = "0987654321-z"
$ws.Range("D6").Value = "Programing"
$ws.Range("E6").Value = 19.0

$ws.Range("A7").Value = "C++"
$ws.Range("B7").Value = "MY_sig"
$ws.Range("C7").Value = "543211234-z"
$ws.Range("D7").Value = "Programing"
$ws.Range("E7").Value = 11.0

$ws.Range("A8").Value = "Python"
$ws.Range("B8").Value = "Aquib"
$ws.Range("C8").Value = "6789009876-z"
$ws.Range("D8").Value = "Programing"
$ws.Range("E8").Value = 19

$ws.Range("A9").Value = "HTML"
$ws.Range("B9").Value = "Devkkumar"
# C9 ("12345678901") is a pure-digit ISBN that must stay TEXT (not auto-convert to a number).
# Writing it through a formula result and then Copy/PasteSpecial(values) preserves the
# shared-string text type without introducing any extra cell style.
$ws.Range("C9").Formula = "=""12345678901"""
$ws.Range("C9").Copy()
$ws.Range("C9").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range("D9").Value = "Programing"
$ws.Range("E9").Value = 21.0

$ws.Range("D8").Select()
